# OpenTBS 1.9.1-beta-2014-07-22 : credit features
#
# The canonical diff for this revision only touches the internal axis
# identifiers of the 3-D bar chart on slide 3 ("Graphique 3"):
#   c:axId   95843456 -> 61990016   (category axis)
#   c:axId   95844992 -> 61991552   (value axis)
#   c:crossAx references updated to match the swap.
#
# These axId/crossAx values are PowerPoint's own internal bookkeeping
# identifiers for linking <c:catAx>/<c:valAx> to the series that use
# them; the Chart/Axis COM object model does not expose a settable
# identifier property for them (there is no Axis.Id / Axis.AxisID on
# the real object model - PowerPoint assigns/rewrites these only when
# it serialises the chart part itself). Nothing about the chart's
# visible configuration (scaling, gridlines, tick marks, crossing
# behaviour, 3-D shape, etc.) changes in this revision, so there is no
# other observable edit to replay here.
#
# Reach the two axis objects the diff is about through the supported
# object model and make a best-effort, non-fatal attempt at renaming
# their identifiers in case the host implements it; this is a no-op
# everywhere the property genuinely isn't settable.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)
$c = $sh.Chart

$catAx = $c.Axes(1)
$valAx = $c.Axes(2)

try { $catAx.AxisID = 61990016 } catch { }
try { $valAx.AxisID = 61991552 } catch { }
